# Update existing rows: refresh the execution-time timestamps for the first
# test case, and turn what used to be row 3 ("Verify Home Page Loads
# Successfully") into a new "User Login with Invalid Credentials" test case
# with its own execution time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("User Login with Valid Credentials") just got re-run -> new timestamp.
$ws.Cells.Item(2, 3).Value = "Thu Mar 27 12:56:16 IST 2025"

# Row 3 used to be "Verify Home Page Loads Successfully"; it is now repurposed
# as the new "User Login with Invalid Credentials" test case, executed at a
# fresh timestamp. Status/comment columns are unchanged.
$ws.Cells.Item(3, 1).Value = "User Login with Invalid Credentials"
$ws.Cells.Item(3, 3).Value = "Thu Mar 27 12:56:21 IST 2025"

# Append a new row 4 that carries the original "Verify Home Page Loads
# Successfully" test case (re-run at the same new timestamp as row 3).
$ws.Cells.Item(4, 1).Value = "Verify Home Page Loads Successfully"
$ws.Cells.Item(4, 2).Value = "PASSED"
$ws.Cells.Item(4, 3).Value = "Thu Mar 27 12:56:21 IST 2025"
$ws.Cells.Item(4, 4).Value = "Test executed successfully."
